# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to store the literal text (avoids Excel auto-converting
    # numeric-looking strings like '506.80' or '1.00' into numbers, which would
    # drop significant trailing zeros / trailing-decimal formatting).
    $range.NumberFormat = '@'
    $range.Value = $text
    $range.Style = 'Normal'
}

Set-TextValue $ws.Range('D2') '56.335.72'
$ws.Range('E2').Value = '  -1.17%  '
Set-TextValue $ws.Range('D3') '2.372.46'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue $ws.Range('D5') '506.80'
$ws.Range('E5').Value = '  +0.37%  '
Set-TextValue $ws.Range('D6') '130.24'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue $ws.Range('D8') '0.545'
$ws.Range('E8').Value = '  -1.44%  '
Set-TextValue $ws.Range('D9') '2.379.89'
$ws.Range('E9').Value = '  -1.15%  '
Set-TextValue $ws.Range('D10') '0.0987'
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('E11').Value = '  -0.10%  '
Set-TextValue $ws.Range('D12') '4.89'
$ws.Range('E12').Value = '  +7.19%  '
$ws.Range('E13').Value = '  +2.21%  '
Set-TextValue $ws.Range('D14') '2.793.72'
$ws.Range('E14').Value = '  -1.15%  '
Set-TextValue $ws.Range('D15') '56.295.55'
$ws.Range('E15').Value = '  -1.12%  '
Set-TextValue $ws.Range('D16') '21.55'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('E17').Value = '  -0.51%  '
Set-TextValue $ws.Range('D18') '2.348.39'
$ws.Range('E18').Value = '  -3.74%  '
Set-TextValue $ws.Range('D19') '10.04'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('E20').Value = '  +0.29%  '
Set-TextValue $ws.Range('D21') '309.68'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  -0.48%  '
Set-TextValue $ws.Range('D23') '1.00'
$ws.Range('E23').Value = '  +0.14%  '
Set-TextValue $ws.Range('D24') '65.83'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('E25').Value = '  -0.02%  '
Set-TextValue $ws.Range('D26') '0.372'
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('E27').Value = '  -2.47%  '
Set-TextValue $ws.Range('D28') '7.22'
$ws.Range('E28').Value = '  -2.91%  '
Set-TextValue $ws.Range('D29') '172.58'
$ws.Range('E29').Value = '  -0.97%  '
Set-TextValue $ws.Range('D30') '0.0₃0712'
$ws.Range('E30').Value = '  -1.55%  '
Set-TextValue $ws.Range('D31') '1.65'
$ws.Range('E31').Value = '  -1.72%  '
Set-TextValue $ws.Range('D32') '5.86'
$ws.Range('E32').Value = '  -1.34%  '
Set-TextValue $ws.Range('D34') '0.995'
$ws.Range('E34').Value = '  -0.08%  '
Set-TextValue $ws.Range('D35') '1.07'
$ws.Range('E35').Value = '  -4.15%  '
Set-TextValue $ws.Range('D36') '17.62'
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E37').Value = '  -1.15%  '
Set-TextValue $ws.Range('D38') '3.70'
$ws.Range('E38').Value = '  -3.36%  '
Set-TextValue $ws.Range('D39') '0.830'
$ws.Range('E39').Value = '  +1.88%  '
Set-TextValue $ws.Range('D40') '36.38'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  -3.58%  '
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D43') '125.79'
$ws.Range('E43').Value = '  -5.43%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D44') '4.73'
$ws.Range('E44').Value = '  -4.82%  '
$ws.Range('E45').Value = '  -0.63%  '
Set-TextValue $ws.Range('D46') '0.0899'
$ws.Range('E46').Value = '  -1.21%  '
Set-TextValue $ws.Range('D47') '237.36'
$ws.Range('E47').Value = '  -5.89%  '
Set-TextValue $ws.Range('D48') '0.0482'
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('E49').Value = '  -1.40%  '
Set-TextValue $ws.Range('D50') '17.00'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('E51').Value = '  +0.04%  '
